$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 731 (the "「ガザ地区」" post). This shifts all subsequent rows up by one,
# matching the author's commit that removed this post row.
$ws.Rows("731:731").Delete()
